# Append three new day-log entries to the tracker sheet (rows 17-19),
# matching the "actually started implementing the game - ice puzzle" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - 2024-12-09 (serial 45635)
$ws.Range("A17").Value = 45635
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("B17").Value = "started an actual boss fight and changes some ideas in the plan"
$ws.Range("C17").Value = 7

# Row 18 - 2024-12-10 (serial 45636)
$ws.Range("A18").Value = 45636
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("B18").Value = "finished most of that boss fight and did some spritework on the arena"
$ws.Range("C18").Value = 6

# Row 19 - 2024-12-11 (serial 45637)
$ws.Range("A19").Value = 45637
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("B19").Value = "animations for the bossfight"
$ws.Range("C19").Value = 6

# Move the active selection to the last entered cell, like Excel would
# leave it after typing the final value.
$ws.Range("C19").Select()
